# Commit: "Add stuff to DJ-02"
#
# The "Demo" slide (title "Demo", with the DJ-03-Model-Single.txt link) that
# sits at slide position 20 is removed from the deck. All following slides
# ("Summary", "Acknowledgements / Contributions") shift up by one position.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)
$s.Delete()
